# Remove the "numberOfPages" column (F) from Sheet1.
# This shifts the former "documentType" column (G) left into column F,
# matching the committed change (column F's numberOfPages/1234/a data is
# gone, and documentType/BAD/BOOK/THESIS now occupies column F).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns("F").Select() | Out-Null
$ws.Columns("F").Delete() | Out-Null
